# "got lat longs for Elahi" -----------------------------------------------
# Adds four new columns (hadLong, hadLat, hadYear1, hadYear2) to the
# study_table sheet right before the existing "methods" column, fills them
# in for all three studies, and also fills in the previously-missing
# study_long / study_lat values for the Elahi2016 row (row 4).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("study_table")

# ---------------------------------------------------------------------
# 1. Fill in the missing study_long / study_lat for Elahi2016 (row 4)
#    *before* inserting columns, so the new cells that appear to its
#    right on that row pick up the correct (non-default) cell style.
# ---------------------------------------------------------------------
$ws.Range("D4").Value = -121.9049
$ws.Range("E4").Value = 36.62186

# ---------------------------------------------------------------------
# 2. Insert four new blank columns at F:I (old columns F.. shift to J..)
# ---------------------------------------------------------------------
$ws.Range("F1:I1").EntireColumn.Insert()

# ---------------------------------------------------------------------
# 3. New header labels for the inserted columns
# ---------------------------------------------------------------------
$ws.Range("F1").Value = "hadLong"
$ws.Range("G1").Value = "hadLat"
$ws.Range("H1").Value = "hadYear1"
$ws.Range("I1").Value = "hadYear2"

# ---------------------------------------------------------------------
# 4. Fill in the new columns, row by row
# ---------------------------------------------------------------------
# Row 2 - Roy2003
$ws.Range("F2").Value = -117.824
$ws.Range("G2").Value = 33.542
$ws.Range("H2").Formula = "=1869 - 10"
$ws.Range("I2").Value = 2001

# Row 3 - Fisher2009
$ws.Range("F3").Value = -68.385
$ws.Range("G3").Value = 44.234
$ws.Range("H3").Formula = "=1915-10"
$ws.Range("I3").Value = 2007

# Row 4 - Elahi2016
$ws.Range("F4").Value = -121.895
$ws.Range("G4").Value = 36.629
$ws.Range("H4").Formula = "= 1947 - 10"
$ws.Range("I4").Value = 2015

# ---------------------------------------------------------------------
# 5. Column widths: new F:I columns match the old "study_lat" width.
#    The old F..J columns (now shifted to J..N) already keep their
#    original widths automatically after EntireColumn.Insert(), so they
#    are intentionally left untouched here.
# ---------------------------------------------------------------------
$ws.Range("F1:I1").EntireColumn.ColumnWidth = 9

# ---------------------------------------------------------------------
# 6. Update the active selection to I3, as in the authored workbook
# ---------------------------------------------------------------------
$ws.Range("I3").Select()

# ---------------------------------------------------------------------
# 7. Move / resize the workbook window to match the authored file
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.Left = 4480
$win.Top = 3100
$win.Width = 21460
$win.Height = 9440
